$wb = $excel.ActiveWorkbook

# New "想去人数" (want-to-go count) values for rows 2-7, column F
$newValues = @{
    2 = 20
    3 = 1814
    4 = 551
    5 = 1167
    6 = 6095
    7 = 142
}

# Both "展览" and "全部类型" sheets carry the same data table and need updating
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $newValues.Keys) {
        $ws.Cells.Item($row, 6).Value = $newValues[$row]
    }
}
